# Adds two new "Logo (not exists ...)" sample paragraphs (with bookmarks +
# an inline image each) right after the existing "Logo (ratio size)"
# paragraph, switches that paragraph's mark formatting from color+lang to
# noProof, and appends one more blank paragraph (color+lang) before the
# document's trailing blank paragraph.

$d = $word.ActiveDocument

# Common namespace declarations used by every InsertXML payload below.
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' + `
      'xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" ' + `
      'xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" ' + `
      'xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" ' + `
      'xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"'

function New-PackageXml($bodyXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
           '<pkg:xmlData><w:document ' + $ns + '><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

function New-ImagePara($bookmarkId, $bookmarkName, $docPrId, $text) {
    return '<w:p><w:pPr><w:rPr><w:color w:val="FF0000"/><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
        '<w:r><w:rPr><w:noProof/></w:rPr><w:t xml:space="preserve">' + $text + '</w:t></w:r>' + `
        '<w:bookmarkStart w:id="' + $bookmarkId + '" w:name="' + $bookmarkName + '"/>' + `
        '<w:r><w:rPr><w:noProof/></w:rPr><w:drawing>' + `
        '<wp:inline distT="0" distB="0" distL="0" distR="0">' + `
        '<wp:extent cx="266700" cy="285750"/>' + `
        '<wp:effectExtent l="19050" t="0" r="0" b="0"/>' + `
        '<wp:docPr id="' + $docPrId + '" name="Image 0" descr="template.png"/>' + `
        '<wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr>' + `
        '<a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture">' + `
        '<pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture">' + `
        '<pic:nvPicPr><pic:cNvPr id="0" name="template.png"/><pic:cNvPicPr/></pic:nvPicPr>' + `
        '<pic:blipFill><a:blip r:embed="rId5" cstate="print"/><a:stretch><a:fillRect/></a:stretch></pic:blipFill>' + `
        '<pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="266700" cy="285750"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom></pic:spPr>' + `
        '</pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r>' + `
        '<w:bookmarkEnd w:id="' + $bookmarkId + '"/></w:p>'
}

# 1) Locate the "Logo (ratio size) : " paragraph (it owns the ratioSizeLogo
#    bookmark) and swap its paragraph-mark run properties from
#    color+lang to a plain noProof.
$ratioBookmark = $d.Bookmarks.Item("ratioSizeLogo")
$ratioParaIndex = $d.Range(0, $ratioBookmark.Range.Start).Paragraphs.Count
$pRatio = $d.Paragraphs.Item($ratioParaIndex)
$rRatio = $pRatio.Range
$xmlRatio = New-PackageXml('<w:p><w:pPr><w:rPr><w:noProof/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:noProof/></w:rPr><w:t xml:space="preserve">Logo (ratio size) : </w:t></w:r>' + `
    '<w:bookmarkStart w:id="3" w:name="ratioSizeLogo"/>' + `
    '<w:r><w:rPr><w:noProof/></w:rPr><w:drawing>' + `
    '<wp:inline distT="0" distB="0" distL="0" distR="0">' + `
    '<wp:extent cx="266700" cy="285750"/>' + `
    '<wp:effectExtent l="19050" t="0" r="0" b="0"/>' + `
    '<wp:docPr id="4" name="Image 0" descr="template.png"/>' + `
    '<wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr>' + `
    '<a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture">' + `
    '<pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture">' + `
    '<pic:nvPicPr><pic:cNvPr id="0" name="template.png"/><pic:cNvPicPr/></pic:nvPicPr>' + `
    '<pic:blipFill><a:blip r:embed="rId5" cstate="print"/><a:stretch><a:fillRect/></a:stretch></pic:blipFill>' + `
    '<pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="266700" cy="285750"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom></pic:spPr>' + `
    '</pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r>' + `
    '<w:bookmarkEnd w:id="3"/></w:p>')
$rRatio.InsertXML($xmlRatio)

# 2) Insert the "Logo (not exists + remove « template » image)" paragraph
#    right after it (bookmark id 4, image docPr id 5).
$pRatio2 = $d.Paragraphs.Item($ratioParaIndex)
$pRatio2.Range.InsertParagraphAfter()
$pRemove = $d.Paragraphs.Item($ratioParaIndex + 1)
$pRemove.Range.InsertXML((New-PackageXml (New-ImagePara 4 "imageNotExistsAndRemoveImageTemplate" 5 "Logo (not exists + remove « template » image) : ")))

# 3) Insert the "Logo (not exists + keep « template » image)" paragraph
#    after that one (bookmark id 5, image docPr id 6).
$pRemove2 = $d.Paragraphs.Item($ratioParaIndex + 1)
$pRemove2.Range.InsertParagraphAfter()
$pKeep = $d.Paragraphs.Item($ratioParaIndex + 2)
$pKeep.Range.InsertXML((New-PackageXml (New-ImagePara 5 "imageNotExistsAndKeepImageTemplate" 6 "Logo (not exists + keep « template » image) : ")))

# 4) Append one more empty paragraph (color+lang mark formatting) before
#    the document's original trailing empty paragraph.
$pKeep2 = $d.Paragraphs.Item($ratioParaIndex + 2)
$pKeep2.Range.InsertParagraphAfter()
$pBlank = $d.Paragraphs.Item($ratioParaIndex + 3)
$xmlBlank = New-PackageXml('<w:p><w:pPr><w:rPr><w:color w:val="FF0000"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>')
$pBlank.Range.InsertXML($xmlBlank)
